# Add "2022-Q1" sheet (new quarter data) right before the "总计" (totals) sheet,
# and update the "总计" sheet with a new first data row for 2022-Q1,
# shifting the previous rows down.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" worksheet just before the totals sheet ---
$totalsSheetBefore = $wb.Worksheets.Item(3)
$formatSource = $wb.Worksheets.Item(2)   # "2021-Q4" sheet - used as a formatting template

$newSheet = $wb.Worksheets.Add($totalsSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet references resolve by live index in this environment, so after
# inserting a sheet "before" position 3, the object that used to be at index 3
# ("总计") has shifted to index 4. Re-fetch it by its new index.
$totalsSheet = $wb.Worksheets.Item(4)

# Copy header-row formatting (bold + border + centered) from the template sheet
$formatSource.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy column-A data formatting (bold + border + centered) from the template sheet
$formatSource.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# --- Header row ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data row ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'002271"
$newSheet.Range("C2").Value = "招商安弘灵活配置混合"
$newSheet.Range("D2").Value = "'0.50"
$newSheet.Range("E2").Value = "'72.34"
$newSheet.Range("F2").Value = "'3.28"
$newSheet.Range("G2").Value = "'0.0164"
$newSheet.Range("H2").Value = 7

# --- 2. Update the "总计" (totals) sheet: add a 2022-Q1 row on top, push others down ---
# Row 4 is brand new (the sheet previously only had rows 1-3), so it needs the
# same "column A" formatting (bold + border + centered) as the existing rows.
$totalsSheet.Range("A2").Copy()
$totalsSheet.Range("A4").PasteSpecial(-4122)

$totalsSheet.Range("A4").Value = 2
$totalsSheet.Range("B4").Value = "2020-Q4"
$totalsSheet.Range("C4").Value = 3
$totalsSheet.Range("D4").Value = 0.21

$totalsSheet.Range("A3").Value = 1
$totalsSheet.Range("B3").Value = "2021-Q4"
$totalsSheet.Range("C3").Value = 4
$totalsSheet.Range("D3").Value = 1.08

$totalsSheet.Range("A2").Value = 0
$totalsSheet.Range("B2").Value = "2022-Q1"
$totalsSheet.Range("C2").Value = 1
$totalsSheet.Range("D2").Value = 0.02
